$wb = $excel.ActiveWorkbook

# --- Sheet1: delete the now-unwanted tail rows (45-87), which only ever
#     held a running index in column A, and shrink the sheet accordingly ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A45:A87").EntireRow.Delete()

# --- Make Sheet1 the active / displayed sheet (instead of Sheet3), with
#     the view scrolled down near the bottom of the (now smaller) data and
#     cell F65 selected ---
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 45
$win.ScrollColumn = 1
$ws1.Range("F65").Select()

Write-Output "done"
